$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage (matching original inline-string cell type) without altering
# the cells visible style: temporarily apply a Text number format while writing
# the value, then restore the cell to the default "Normal" style.
$cells = @(
    @{Addr = "D2"; Val = "27.978.99"}
    @{Addr = "E2"; Val = "  -0.23%  "}
    @{Addr = "D3"; Val = "1.911.46"}
    @{Addr = "E3"; Val = "  -0.06%  "}
    @{Addr = "D4"; Val = "0.9972"}
    @{Addr = "E4"; Val = "  -0.84%  "}
    @{Addr = "D5"; Val = "313.55"}
    @{Addr = "E5"; Val = "  -0.71%  "}
    @{Addr = "D6"; Val = "0.9982"}
    @{Addr = "E6"; Val = "  -0.66%  "}
    @{Addr = "D7"; Val = "0.4999"}
    @{Addr = "E7"; Val = "  +3.60%  "}
    @{Addr = "D8"; Val = "0.3823"}
    @{Addr = "E8"; Val = "  +0.24%  "}
    @{Addr = "D9"; Val = "0.07308"}
    @{Addr = "E9"; Val = "  -0.78%  "}
    @{Addr = "D10"; Val = "0.9139"}
    @{Addr = "E10"; Val = "  -2.20%  "}
    @{Addr = "D11"; Val = "21.26"}
    @{Addr = "E11"; Val = "  +2.05%  "}
    @{Addr = "D12"; Val = "0.07686"}
    @{Addr = "E12"; Val = "  -1.29%  "}
    @{Addr = "D13"; Val = "1.892.86"}
    @{Addr = "E13"; Val = "  -1.18%  "}
    @{Addr = "D14"; Val = "5.517"}
    @{Addr = "E14"; Val = "  +0.25%  "}
    @{Addr = "D15"; Val = "92.74"}
    @{Addr = "E15"; Val = "  +0.68%  "}
    @{Addr = "D16"; Val = "0.9978"}
    @{Addr = "E16"; Val = "  -0.79%  "}
    @{Addr = "D17"; Val = "0.000008751"}
    @{Addr = "E17"; Val = "  -1.35%  "}
    @{Addr = "D18"; Val = "0.9975"}
    @{Addr = "E18"; Val = "  -0.72%  "}
    @{Addr = "D19"; Val = "27.975.45"}
    @{Addr = "E19"; Val = "  -0.36%  "}
    @{Addr = "D20"; Val = "14.67"}
    @{Addr = "E20"; Val = "  -0.79%  "}
    @{Addr = "D21"; Val = "5.186"}
    @{Addr = "E21"; Val = "  +0.28%  "}
    @{Addr = "E22"; Val = "  -0.58%  "}
    @{Addr = "D23"; Val = "6.602"}
    @{Addr = "E23"; Val = "  -0.61%  "}
    @{Addr = "E24"; Val = "  -2.00%  "}
    @{Addr = "D25"; Val = "1.847"}
    @{Addr = "E25"; Val = "  -3.56%  "}
    @{Addr = "D26"; Val = "2.224"}
    @{Addr = "E26"; Val = "  +4.34%  "}
    @{Addr = "D27"; Val = "18.43"}
    @{Addr = "E27"; Val = "  -0.34%  "}
    @{Addr = "D28"; Val = "115.43"}
    @{Addr = "E28"; Val = "  -1.38%  "}
    @{Addr = "D29"; Val = "4.916"}
    @{Addr = "E29"; Val = "  -1.11%  "}
    @{Addr = "D30"; Val = "0.09027"}
    @{Addr = "E30"; Val = "  +0.75%  "}
    @{Addr = "D31"; Val = "3.206"}
    @{Addr = "E31"; Val = "  -3.08%  "}
    @{Addr = "D32"; Val = "4.862"}
    @{Addr = "E32"; Val = "  +3.99%  "}
    @{Addr = "D33"; Val = "1.236"}
    @{Addr = "E33"; Val = "  -2.39%  "}
    @{Addr = "D34"; Val = "0.7737"}
    @{Addr = "E34"; Val = "  -0.82%  "}
    @{Addr = "D35"; Val = "0.02092"}
    @{Addr = "E35"; Val = "  +1.77%  "}
    @{Addr = "D36"; Val = "2.573"}
    @{Addr = "E36"; Val = "  -1.62%  "}
    @{Addr = "D37"; Val = "3.066"}
    @{Addr = "E37"; Val = "  +2.36%  "}
    @{Addr = "E38"; Val = "  -1.67%  "}
    @{Addr = "D39"; Val = "0.5573"}
    @{Addr = "E39"; Val = "  +1.11%  "}
    @{Addr = "D40"; Val = "0.05287"}
    @{Addr = "E40"; Val = "  -0.61%  "}
    @{Addr = "D41"; Val = "6.895"}
    @{Addr = "E41"; Val = "  -2.05%  "}
    @{Addr = "D42"; Val = "8.526"}
    @{Addr = "E42"; Val = "  +0.41%  "}
    @{Addr = "D43"; Val = "112.88"}
    @{Addr = "E43"; Val = "  +3.96%  "}
    @{Addr = "D44"; Val = "0.1523"}
    @{Addr = "E44"; Val = "  -0.41%  "}
    @{Addr = "D45"; Val = "10.63"}
    @{Addr = "E45"; Val = "  -0.91%  "}
    @{Addr = "D46"; Val = "0.4845"}
    @{Addr = "E46"; Val = "  +0.34%  "}
    @{Addr = "D47"; Val = "0.9986"}
    @{Addr = "E47"; Val = "  -0.63%  "}
    @{Addr = "D48"; Val = "1.640"}
    @{Addr = "E48"; Val = "  -0.57%  "}
    @{Addr = "D49"; Val = "67.56"}
    @{Addr = "E49"; Val = "  -0.94%  "}
    @{Addr = "D50"; Val = "0.06046"}
    @{Addr = "E50"; Val = "  -0.58%  "}
    @{Addr = "D51"; Val = "0.9073"}
)

foreach ($item in $cells) {
    $rng = $ws.Range($item.Addr)
    $rng.NumberFormat = "@"
    $rng.Value = $item.Val
    $rng.Style = "Normal"
}
